# Análise funcionando com filtro por Estratégia
#
# - Remove the stray "ads"/"adsa" scratch row (row 6) from the "params"
#   sheet, shifting the trailing styled-but-empty cell up to B6.
# - Update the remembered selection on "params" to B9 and on "levers" to
#   D3, while keeping "params" as the active sheet/tab.

$wb = $excel.ActiveWorkbook

$wsParams = $wb.Worksheets.Item("params")
$wsLevers = $wb.Worksheets.Item("levers")

# Set the levers sheet selection first (without leaving it as the active
# tab once we're done).
[void]$wsLevers.Range("D3").Select()

# Make params the active sheet again, then perform its edits/selection
# last so it ends up as the tab that is actually active when saved.
$wsParams.Activate()
$wsParams.Rows("6").Delete()
[void]$wsParams.Range("B9").Select()
